# Applies the diff: inserts one new weekly data row before row 233
# (shifting the previous rows 233-237 down to 234-238) and fills the
# new row 233 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 233, pushing existing rows 233-237 down to 234-238.
$ws.Rows.Item(233).EntireRow.Insert()

# Populate the newly inserted row 233 with the new data record.
$ws.Range("A233").Value = 2
$ws.Range("B233").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C233").Value = "Coquimbo"
$ws.Range("D233").Value = 45267
$ws.Range("E233").Value = 4
$ws.Range("F233").Value = "Fruta"
$ws.Range("G233").Value = 100109
$ws.Range("H233").Value = "Uva"
$ws.Range("I233").Value = 100109001
$ws.Range("J233").Value = "Uva"
$ws.Range("K233").Value = "Flame Seedless"
$ws.Range("L233").Value = "Primera"
$ws.Range("M233").Value = 1100
$ws.Range("N233").Value = 13000
$ws.Range("O233").Value = 14000
$ws.Range("P233").Value = 13500
$ws.Range("Q233").Value = "`$/bandeja 10 kilos"
$ws.Range("R233").Value = "Provincia del Elquí"
$ws.Range("S233").Value = 1350
$ws.Range("T233").Value = 10
